# Auto update Excel log 2026-02-04 14:30:21
# Appends newly-captured sensor readings to the PIR, Humidity and
# Temperature logs. Each row: Date, Timestamp, Hour, Location, Value, Status

$wb = $excel.ActiveWorkbook

function Append-Rows {
    param($SheetName, $StartRow, $EndRow, $Rows)

    $ws = $wb.Worksheets.Item($SheetName)

    # Column A holds plain "yyyy-mm-dd" log dates that must stay literal
    # text (matching every pre-existing row) instead of being auto-parsed
    # into a date serial number, so force Text format before writing them.
    $ws.Range("A" + $StartRow + ":A" + $EndRow).NumberFormat = "@"

    $r = $StartRow
    foreach ($row in $Rows) {
        $ws.Cells.Item($r, 1).Value = $row[0]
        $ws.Cells.Item($r, 2).Value = $row[1]
        $ws.Cells.Item($r, 3).Value = $row[2]
        $ws.Cells.Item($r, 4).Value = $row[3]
        $ws.Cells.Item($r, 5).Value = $row[4]
        $ws.Cells.Item($r, 6).Value = $row[5]
        $r = $r + 1
    }
}

# ---- PIR sheet: rows 338-349 ----
$pirRows = @(
    ,@("2026-02-04","14:29:15","14:00","Bathroom","No Motion","Inactive")
    ,@("2026-02-04","14:29:17","14:00","Bathroom","No Motion","Inactive")
    ,@("2026-02-04","14:29:21","14:00","Bathroom","No Motion","Inactive")
    ,@("2026-02-04","14:29:25","14:00","Bathroom","Motion Detected","Active")
    ,@("2026-02-04","14:29:34","14:00","Bathroom","No Motion","Inactive")
    ,@("2026-02-04","14:29:39","14:00","Bathroom","No Motion","Inactive")
    ,@("2026-02-04","14:29:43","14:00","Bathroom","Motion Detected","Active")
    ,@("2026-02-04","14:29:50","14:00","Bathroom","No Motion","Inactive")
    ,@("2026-02-04","14:29:54","14:00","Bathroom","Motion Detected","Active")
    ,@("2026-02-04","14:30:10","14:00","Bathroom","No Motion","Inactive")
    ,@("2026-02-04","14:30:11","14:00","Bathroom","Motion Detected","Active")
    ,@("2026-02-04","14:30:14","14:00","Bathroom","No Motion","Inactive")
)
Append-Rows "PIR" 338 349 $pirRows

# ---- Humidity sheet: rows 278-287 ----
$humidityRows = @(
    ,@("2026-02-04","14:29:16","14:00","Bathroom","79.7%","Active")
    ,@("2026-02-04","14:29:18","14:00","Bathroom","78.8%","Active")
    ,@("2026-02-04","14:29:21","14:00","Bathroom","79.7%","Active")
    ,@("2026-02-04","14:29:26","14:00","Bathroom","78.8%","Active")
    ,@("2026-02-04","14:29:32","14:00","Bathroom","79.7%","Active")
    ,@("2026-02-04","14:29:37","14:00","Bathroom","78.7%","Active")
    ,@("2026-02-04","14:29:42","14:00","Bathroom","79.7%","Active")
    ,@("2026-02-04","14:29:47","14:00","Bathroom","78.7%","Active")
    ,@("2026-02-04","14:30:08","14:00","Bathroom","79.6%","Active")
    ,@("2026-02-04","14:30:12","14:00","Bathroom","78.6%","Active")
)
$hws = $wb.Worksheets.Item("Humidity")
$hws.Range("E278:E287").NumberFormat = "@"
Append-Rows "Humidity" 278 287 $humidityRows

# ---- Temperature sheet: rows 277-287 ----
$temperatureRows = @(
    ,@("2026-02-04","14:29:14","14:00","Bathroom","24.3C","Active")
    ,@("2026-02-04","14:29:17","14:00","Bathroom","24.3C","Active")
    ,@("2026-02-04","14:29:19","14:00","Bathroom","24.3C","Active")
    ,@("2026-02-04","14:29:22","14:00","Bathroom","24.3C","Active")
    ,@("2026-02-04","14:29:27","14:00","Bathroom","24.3C","Active")
    ,@("2026-02-04","14:29:33","14:00","Bathroom","24.3C","Active")
    ,@("2026-02-04","14:29:38","14:00","Bathroom","24.3C","Active")
    ,@("2026-02-04","14:29:42","14:00","Bathroom","24.3C","Active")
    ,@("2026-02-04","14:29:48","14:00","Bathroom","24.3C","Active")
    ,@("2026-02-04","14:30:09","14:00","Bathroom","24.3C","Active")
    ,@("2026-02-04","14:30:13","14:00","Bathroom","24.4C","Active")
)
Append-Rows "Temperature" 277 287 $temperatureRows
